# DataEngine2.xlsx regression-suite update:
#  - the whole suite now runs clean in one pass except JiraId-1206, so the
#    stale "PASS" results for the first three Test_Steps rows are cleared
#  - row 15's highlight flag (B15) is no longer needed, so its fill is
#    reset back to the normal row style
#  - the Warsaw address postal code used by validate_autoComplete was
#    wrong (00-401) and is corrected to 02-326; the CONCATENATE helper in
#    D4 recalculates on its own

$wb = $excel.ActiveWorkbook

$steps = $wb.Worksheets.Item("Test_Steps")

# Test results have gone stale now that the suite runs end-to-end again -
# clear the cached "PASS" markers in column F (Result) for rows 2-4.
$steps.Range("F2").Value = ""
$steps.Range("F3").Value = ""
$steps.Range("F4").Value = ""

# B15 was highlighted (fill color); drop back to the plain body style by
# copying the format from a neighbouring cell that already has it (C15).
$steps.Range("C15").Copy()
$steps.Range("B15").PasteSpecial(-4122)

# Fix the Warsaw postal code used in the autocomplete validation table.
$autoComplete = $wb.Worksheets.Item("validate_autoComplete")
$autoComplete.Range("B4").Value = "02-326"
